$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (Pull Inventory Report): remove the CU marks in I12 and J12
$ws.Range("I12").ClearContents()
$ws.Range("J12").ClearContents()

# Row 13 (Pull-Inventory Report Definition): remove the CU marks in I13 and J13
$ws.Range("I13").ClearContents()
$ws.Range("J13").ClearContents()

# Row 20 (Pull-Inventory Form): move the CU mark from I20 to J20
$ws.Range("I20").Cut($ws.Range("J20"))

# Row 28 (Pull-Inventory-Line Subform): remove the CU mark in J28 (keep I28)
$ws.Range("J28").ClearContents()

# Update the active selection to match the edited cell
$ws.Range("I20").Select()
